$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column keeps its text formatting so values like "181.70" are not
# auto-converted to numbers by Excel when re-assigned.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '65.535.73'
$ws.Range("E2").Value = '  +1.11%  '
$ws.Range("D3").Value = '3.367.39'
$ws.Range("E3").Value = '  +0.74%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '181.70'
$ws.Range("E5").Value = '  -0.39%  '
$ws.Range("D6").Value = '537.95'
$ws.Range("E6").Value = '  +0.55%  '
$ws.Range("D7").Value = '0.604'
$ws.Range("E7").Value = '  -1.02%  '
$ws.Range("D8").Value = '3.359.81'
$ws.Range("E8").Value = '  +0.52%  '
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("E10").Value = '  +1.19%  '
$ws.Range("D11").Value = '54.94'
$ws.Range("E11").Value = '  -7.01%  '
$ws.Range("E12").Value = '  +5.93%  '
$ws.Range("D13").Value = '0.0000266'
$ws.Range("E13").Value = '  +1.17%  '
$ws.Range("D14").Value = '9.22'
$ws.Range("E14").Value = '  +0.08%  '
$ws.Range("D15").Value = '3.910.67'
$ws.Range("E15").Value = '  +1.08%  '
$ws.Range("E16").Value = '  +1.56%  '
$ws.Range("D17").Value = '3.370.17'
$ws.Range("E17").Value = '  +0.90%  '
$ws.Range("D18").Value = '17.98'
$ws.Range("E18").Value = '  +1.63%  '
$ws.Range("D19").Value = '65.745.79'
$ws.Range("E19").Value = '  +1.60%  '
$ws.Range("D20").Value = '11.48'
$ws.Range("E20").Value = '  +1.88%  '
$ws.Range("E21").Value = '  +0.93%  '
$ws.Range("D22").Value = '392.96'
$ws.Range("E22").Value = '  +3.76%  '
$ws.Range("D23").Value = '11.73'
$ws.Range("E23").Value = '  +3.62%  '
$ws.Range("D24").Value = '4.24'
$ws.Range("E24").Value = '  +7.13%  '
$ws.Range("D25").Value = '83.21'
$ws.Range("E25").Value = '  +2.17%  '
$ws.Range("E26").Value = '  -1.49%  '
$ws.Range("B27").Value = 'LEO'
$ws.Range("C27").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D27").Value = '6.12'
$ws.Range("E27").Value = '  +0.43%  '
$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").Value = '2.83'
$ws.Range("E28").Value = '  +4.52%  '
$ws.Range("D29").Value = '11.56'
$ws.Range("E29").Value = '  -0.19%  '
$ws.Range("D30").Value = '8.45'
$ws.Range("E30").Value = '  -0.29%  '
$ws.Range("D31").Value = '29.51'
$ws.Range("E31").Value = '  +0.76%  '
$ws.Range("D32").Value = '663.60'
$ws.Range("E32").Value = '  +0.56%  '
$ws.Range("D33").Value = '6.75'
$ws.Range("E33").Value = '  -0.47%  '
$ws.Range("E34").Value = '  +0.39%  '
$ws.Range("E35").Value = '  +0.76%  '
$ws.Range("D36").Value = '57.95'
$ws.Range("E36").Value = '  -3.19%  '
$ws.Range("D37").Value = '37.58'
$ws.Range("E37").Value = '  +0.87%  '
$ws.Range("E38").Value = '  +0.06%  '
$ws.Range("E39").Value = '  -0.18%  '
$ws.Range("D40").Value = '0.0₃0777'
$ws.Range("E40").Value = '  +8.53%  '
$ws.Range("D41").Value = '2.75'
$ws.Range("E41").Value = '  +8.62%  '
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  +0.18%  '
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").Value = '3.27'
$ws.Range("E43").Value = '  +16.58%  '
$ws.Range("E44").Value = '  +1.64%  '
$ws.Range("D45").Value = '3.010.41'
$ws.Range("E45").Value = '  +2.10%  '
$ws.Range("E46").Value = '  +1.15%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = '3.25'
$ws.Range("E47").Value = '  +5.33%  '
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").Value = '0.0411'
$ws.Range("E48").Value = '  +2.02%  '
$ws.Range("E49").Value = '  +2.10%  '
$ws.Range("D50").Value = '8.79'
$ws.Range("E50").Value = '  +9.58%  '
$ws.Range("E51").Value = '  +0.25%  '
